{"js": "// Replace each two-digit multiplication problem's text with its updated version.\n// Each [oldText, newText] pair corresponds to one <w:t> run inside a table cell;\n// all old values are unique within the document, so a simple search+replace is safe.\nconst pairs = [\n  [\"17\u00d792=1564\", \"39\u00d765=2535\"],\n  [\"79\u00d758=4582\", \"58\u00d771=4118\"],\n  [\"29\u00d787=2523\", \"78\u00d766=5148\"],\n  [\"20\u00d733=660\", \"71\u00d748=3408\"],\n  [\"74\u00d749=3626\", \"32\u00d745=1440\"],\n  [\"89\u00d716=1424\", \"58\u00d748=2784\"],\n  [\"22\u00d773=1606\", \"42\u00d774=3108\"],\n  [\"85\u00d727=2295\", \"40\u00d715=600\"],\n  [\"96\u00d721=2016\", \"67\u00d716=1072\"],\n  [\"10\u00d753=530\", \"96\u00d767=6432\"],\n  [\"21\u00d748=1008\", \"27\u00d746=1242\"],\n  [\"86\u00d798=8428\", \"43\u00d740=1720\"],\n  [\"25\u00d753=1325\", \"50\u00d745=2250\"],\n  [\"66\u00d793=6138\", \"63\u00d797=6111\"],\n  [\"49\u00d771=3479\", \"56\u00d7100=5600\"],\n  [\"30\u00d718=540\", \"44\u00d7100=4400\"],\n  [\"68\u00d757=3876\", \"92\u00d789=8188\"],\n  [\"14\u00d768=952\", \"70\u00d750=3500\"],\n  [\"76\u00d793=7068\", \"45\u00d769=3105\"],\n  [\"61\u00d745=2745\", \"25\u00d776=1900\"],\n  [\"48\u00d770=3360\", \"81\u00d731=2511\"],\n  [\"46\u00d734=1564\", \"55\u00d710=550\"],\n  [\"52\u00d779=4108\", \"45\u00d775=3375\"],\n  [\"10\u00d717=170\", \"72\u00d796=6912\"],\n  [\"38\u00d741=1558\", \"84\u00d794=7896\"],\n  [\"56\u00d725=1400\", \"41\u00d795=3895\"],\n  [\"85\u00d789=7565\", \"33\u00d773=2409\"],\n  [\"64\u00d741=2624\", \"64\u00d723=1472\"],\n  [\"32\u00d738=1216\", \"97\u00d754=5238\"],\n  [\"42\u00d721=882\", \"38\u00d793=3534\"],\n  [\"64\u00d753=3392\", \"91\u00d787=7917\"],\n  [\"26\u00d772=1872\", \"95\u00d711=1045\"],\n  [\"78\u00d785=6630\", \"90\u00d713=1170\"],\n  [\"97\u00d767=6499\", \"96\u00d742=4032\"],\n  [\"29\u00d710=290\", \"39\u00d795=3705\"],\n  [\"96\u00d789=8544\", \"35\u00d723=805\"],\n  [\"56\u00d779=4424\", \"43\u00d733=1419\"],\n  [\"80\u00d766=5280\", \"34\u00d716=544\"],\n  [\"40\u00d749=1960\", \"14\u00d790=1260\"],\n  [\"95\u00d746=4370\", \"50\u00d751=2550\"],\n  [\"56\u00d722=1232\", \"80\u00d772=5760\"],\n  [\"57\u00d721=1197\", \"10\u00d788=880\"],\n  [\"71\u00d762=4402\", \"21\u00d798=2058\"],\n  [\"42\u00d789=3738\", \"42\u00d795=3990\"],\n  [\"28\u00d771=1988\", \"54\u00d750=2700\"],\n  [\"39\u00d715=585\", \"70\u00d794=6580\"],\n  [\"47\u00d724=1128\", \"53\u00d728=1484\"],\n  [\"28\u00d792=2576\", \"23\u00d794=2162\"],\n  [\"45\u00d738=1710\", \"73\u00d790=6570\"],\n  [\"68\u00d763=4284\", \"45\u00d787=3915\"],\n  [\"33\u00d760=1980\", \"88\u00d725=2200\"],\n  [\"72\u00d767=4824\", \"14\u00d793=1302\"],\n  [\"38\u00d736=1368\", \"59\u00d771=4189\"],\n  [\"44\u00d722=968\", \"37\u00d742=1554\"],\n  [\"36\u00d722=792\", \"84\u00d759=4956\"],\n  [\"56\u00d774=4144\", \"45\u00d797=4365\"],\n  [\"84\u00d744=3696\", \"14\u00d724=336\"],\n  [\"90\u00d712=1080\", \"40\u00d729=1160\"],\n  [\"39\u00d731=1209\", \"53\u00d754=2862\"],\n  [\"74\u00d726=1924\", \"54\u00d728=1512\"],\n  [\"32\u00d7100=3200\", \"47\u00d783=3901\"],\n  [\"13\u00d719=247\", \"69\u00d767=4623\"],\n  [\"36\u00d730=1080\", \"88\u00d732=2816\"],\n  [\"46\u00d743=1978\", \"62\u00d787=5394\"],\n  [\"27\u00d790=2430\", \"30\u00d777=2310\"],\n  [\"49\u00d7100=4900\", \"66\u00d795=6270\"],\n  [\"21\u00d745=945\", \"17\u00d729=493\"],\n  [\"22\u00d757=1254\", \"84\u00d784=7056\"],\n  [\"56\u00d778=4368\", \"29\u00d775=2175\"],\n  [\"37\u00d738=1406\", \"59\u00d799=5841\"],\n  [\"91\u00d718=1638\", \"40\u00d7100=4000\"],\n  [\"20\u00d758=1160\", \"52\u00d764=3328\"],\n  [\"100\u00d732=3200\", \"42\u00d760=2520\"],\n  [\"36\u00d725=900\", \"91\u00d717=1547\"],\n  [\"17\u00d723=391\", \"11\u00d714=154\"],\n  [\"20\u00d764=1280\", \"61\u00d767=4087\"],\n  [\"64\u00d719=1216\", \"19\u00d792=1748\"],\n  [\"67\u00d766=4422\", \"57\u00d748=2736\"],\n  [\"97\u00d714=1358\", \"89\u00d796=8544\"],\n  [\"55\u00d718=990\", \"32\u00d775=2400\"],\n  [\"73\u00d761=4453\", \"72\u00d711=792\"],\n  [\"81\u00d789=7209\", \"76\u00d757=4332\"],\n  [\"13\u00d728=364\", \"54\u00d767=3618\"],\n  [\"50\u00d797=4850\", \"42\u00d749=2058\"],\n  [\"67\u00d764=4288\", \"94\u00d775=7050\"],\n  [\"91\u00d742=3822\", \"43\u00d743=1849\"],\n  [\"87\u00d780=6960\", \"71\u00d719=1349\"],\n  [\"72\u00d715=1080\", \"31\u00d785=2635\"],\n  [\"62\u00d764=3968\", \"38\u00d734=1292\"],\n  [\"27\u00d791=2457\", \"98\u00d716=1568\"],\n  [\"19\u00d742=798\", \"98\u00d766=6468\"],\n  [\"86\u00d736=3096\", \"81\u00d790=7290\"],\n  [\"39\u00d761=2379\", \"28\u00d796=2688\"],\n  [\"11\u00d799=1089\", \"63\u00d752=3276\"],\n  [\"19\u00d768=1292\", \"56\u00d721=1176\"],\n  [\"83\u00d760=4980\", \"79\u00d792=7268\"],\n  [\"39\u00d779=3081\", \"61\u00d790=5490\"],\n  [\"83\u00d748=3984\", \"69\u00d751=3519\"],\n  [\"61\u00d714=854\", \"45\u00d750=2250\"],\n  [\"65\u00d716=1040\", \"77\u00d795=7315\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication problem's text with its updated version.\n# Each ,@(oldText, newText) pair corresponds to one table-cell run; all old\n# values are unique within the document, so Find/Replace (first match) is safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"17\u00d792=1564\", \"39\u00d765=2535\")\n    ,@(\"79\u00d758=4582\", \"58\u00d771=4118\")\n    ,@(\"29\u00d787=2523\", \"78\u00d766=5148\")\n    ,@(\"20\u00d733=660\", \"71\u00d748=3408\")\n    ,@(\"74\u00d749=3626\", \"32\u00d745=1440\")\n    ,@(\"89\u00d716=1424\", \"58\u00d748=2784\")\n    ,@(\"22\u00d773=1606\", \"42\u00d774=3108\")\n    ,@(\"85\u00d727=2295\", \"40\u00d715=600\")\n    ,@(\"96\u00d721=2016\", \"67\u00d716=1072\")\n    ,@(\"10\u00d753=530\", \"96\u00d767=6432\")\n    ,@(\"21\u00d748=1008\", \"27\u00d746=1242\")\n    ,@(\"86\u00d798=8428\", \"43\u00d740=1720\")\n    ,@(\"25\u00d753=1325\", \"50\u00d745=2250\")\n    ,@(\"66\u00d793=6138\", \"63\u00d797=6111\")\n    ,@(\"49\u00d771=3479\", \"56\u00d7100=5600\")\n    ,@(\"30\u00d718=540\", \"44\u00d7100=4400\")\n    ,@(\"68\u00d757=3876\", \"92\u00d789=8188\")\n    ,@(\"14\u00d768=952\", \"70\u00d750=3500\")\n    ,@(\"76\u00d793=7068\", \"45\u00d769=3105\")\n    ,@(\"61\u00d745=2745\", \"25\u00d776=1900\")\n    ,@(\"48\u00d770=3360\", \"81\u00d731=2511\")\n    ,@(\"46\u00d734=1564\", \"55\u00d710=550\")\n    ,@(\"52\u00d779=4108\", \"45\u00d775=3375\")\n    ,@(\"10\u00d717=170\", \"72\u00d796=6912\")\n    ,@(\"38\u00d741=1558\", \"84\u00d794=7896\")\n    ,@(\"56\u00d725=1400\", \"41\u00d795=3895\")\n    ,@(\"85\u00d789=7565\", \"33\u00d773=2409\")\n    ,@(\"64\u00d741=2624\", \"64\u00d723=1472\")\n    ,@(\"32\u00d738=1216\", \"97\u00d754=5238\")\n    ,@(\"42\u00d721=882\", \"38\u00d793=3534\")\n    ,@(\"64\u00d753=3392\", \"91\u00d787=7917\")\n    ,@(\"26\u00d772=1872\", \"95\u00d711=1045\")\n    ,@(\"78\u00d785=6630\", \"90\u00d713=1170\")\n    ,@(\"97\u00d767=6499\", \"96\u00d742=4032\")\n    ,@(\"29\u00d710=290\", \"39\u00d795=3705\")\n    ,@(\"96\u00d789=8544\", \"35\u00d723=805\")\n    ,@(\"56\u00d779=4424\", \"43\u00d733=1419\")\n    ,@(\"80\u00d766=5280\", \"34\u00d716=544\")\n    ,@(\"40\u00d749=1960\", \"14\u00d790=1260\")\n    ,@(\"95\u00d746=4370\", \"50\u00d751=2550\")\n    ,@(\"56\u00d722=1232\", \"80\u00d772=5760\")\n    ,@(\"57\u00d721=1197\", \"10\u00d788=880\")\n    ,@(\"71\u00d762=4402\", \"21\u00d798=2058\")\n    ,@(\"42\u00d789=3738\", \"42\u00d795=3990\")\n    ,@(\"28\u00d771=1988\", \"54\u00d750=2700\")\n    ,@(\"39\u00d715=585\", \"70\u00d794=6580\")\n    ,@(\"47\u00d724=1128\", \"53\u00d728=1484\")\n    ,@(\"28\u00d792=2576\", \"23\u00d794=2162\")\n    ,@(\"45\u00d738=1710\", \"73\u00d790=6570\")\n    ,@(\"68\u00d763=4284\", \"45\u00d787=3915\")\n    ,@(\"33\u00d760=1980\", \"88\u00d725=2200\")\n    ,@(\"72\u00d767=4824\", \"14\u00d793=1302\")\n    ,@(\"38\u00d736=1368\", \"59\u00d771=4189\")\n    ,@(\"44\u00d722=968\", \"37\u00d742=1554\")\n    ,@(\"36\u00d722=792\", \"84\u00d759=4956\")\n    ,@(\"56\u00d774=4144\", \"45\u00d797=4365\")\n    ,@(\"84\u00d744=3696\", \"14\u00d724=336\")\n    ,@(\"90\u00d712=1080\", \"40\u00d729=1160\")\n    ,@(\"39\u00d731=1209\", \"53\u00d754=2862\")\n    ,@(\"74\u00d726=1924\", \"54\u00d728=1512\")\n    ,@(\"32\u00d7100=3200\", \"47\u00d783=3901\")\n    ,@(\"13\u00d719=247\", \"69\u00d767=4623\")\n    ,@(\"36\u00d730=1080\", \"88\u00d732=2816\")\n    ,@(\"46\u00d743=1978\", \"62\u00d787=5394\")\n    ,@(\"27\u00d790=2430\", \"30\u00d777=2310\")\n    ,@(\"49\u00d7100=4900\", \"66\u00d795=6270\")\n    ,@(\"21\u00d745=945\", \"17\u00d729=493\")\n    ,@(\"22\u00d757=1254\", \"84\u00d784=7056\")\n    ,@(\"56\u00d778=4368\", \"29\u00d775=2175\")\n    ,@(\"37\u00d738=1406\", \"59\u00d799=5841\")\n    ,@(\"91\u00d718=1638\", \"40\u00d7100=4000\")\n    ,@(\"20\u00d758=1160\", \"52\u00d764=3328\")\n    ,@(\"100\u00d732=3200\", \"42\u00d760=2520\")\n    ,@(\"36\u00d725=900\", \"91\u00d717=1547\")\n    ,@(\"17\u00d723=391\", \"11\u00d714=154\")\n    ,@(\"20\u00d764=1280\", \"61\u00d767=4087\")\n    ,@(\"64\u00d719=1216\", \"19\u00d792=1748\")\n    ,@(\"67\u00d766=4422\", \"57\u00d748=2736\")\n    ,@(\"97\u00d714=1358\", \"89\u00d796=8544\")\n    ,@(\"55\u00d718=990\", \"32\u00d775=2400\")\n    ,@(\"73\u00d761=4453\", \"72\u00d711=792\")\n    ,@(\"81\u00d789=7209\", \"76\u00d757=4332\")\n    ,@(\"13\u00d728=364\", \"54\u00d767=3618\")\n    ,@(\"50\u00d797=4850\", \"42\u00d749=2058\")\n    ,@(\"67\u00d764=4288\", \"94\u00d775=7050\")\n    ,@(\"91\u00d742=3822\", \"43\u00d743=1849\")\n    ,@(\"87\u00d780=6960\", \"71\u00d719=1349\")\n    ,@(\"72\u00d715=1080\", \"31\u00d785=2635\")\n    ,@(\"62\u00d764=3968\", \"38\u00d734=1292\")\n    ,@(\"27\u00d791=2457\", \"98\u00d716=1568\")\n    ,@(\"19\u00d742=798\", \"98\u00d766=6468\")\n    ,@(\"86\u00d736=3096\", \"81\u00d790=7290\")\n    ,@(\"39\u00d761=2379\", \"28\u00d796=2688\")\n    ,@(\"11\u00d799=1089\", \"63\u00d752=3276\")\n    ,@(\"19\u00d768=1292\", \"56\u00d721=1176\")\n    ,@(\"83\u00d760=4980\", \"79\u00d792=7268\")\n    ,@(\"39\u00d779=3081\", \"61\u00d790=5490\")\n    ,@(\"83\u00d748=3984\", \"69\u00d751=3519\")\n    ,@(\"61\u00d714=854\", \"45\u00d750=2250\")\n    ,@(\"65\u00d716=1040\", \"77\u00d795=7315\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
